$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update status values (Status column D) and actual-time values (column F) ---
# Row 5: status -> Slutförd, verklig tid 2 -> 7
$ws.Range("D5").Value = "Slutförd"
$ws.Range("F5").Value = 7

# Row 6: status -> Slutförd
$ws.Range("D6").Value = "Slutförd"

# Row 7: status -> Påbörjad, verklig tid 0 -> 1
$ws.Range("D7").Value = "Påbörjad"
$ws.Range("F7").Value = 1

# Row 10: status -> Slutförd, verklig tid 0 -> 2
$ws.Range("D10").Value = "Slutförd"
$ws.Range("F10").Value = 2

# Row 11: verklig tid 0 -> 4
$ws.Range("F11").Value = 4

# --- New row 12: continuation of the requirement list (ID 15) ---
$ws.Range("A12").Value = 15

# --- Move the "Summa" totals row from row 13 down to row 15, extending the sums ---
$ws.Range("D15").Value = "Summa"
$ws.Range("E15").Formula = "=SUM(E5:E13)"
$ws.Range("F15").Formula = "=SUM(F5:F13)"
$ws.Range("D13:F13").ClearContents()

# --- Highlight rows 8 and 9 (ID column) with a solid red fill ---
$ws.Range("A8").Interior.Color = 255
$ws.Range("A9").Interior.Color = 255

# --- Update the active selection shown when the workbook is opened ---
$ws.Range("C13").Select()
